$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, bordered, centered) from E1 to F1 so the new
# "time_taken" header cell reuses the existing header style (s="1").
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("F1").Value = "time_taken"

$ws.Range("F2").Value = "2021-10-05 13:39:53.448640"
$ws.Range("F3").Value = "2021-10-05 13:39:53.448651"
$ws.Range("F4").Value = "2021-10-05 13:39:53.448655"
$ws.Range("F5").Value = "2021-10-05 13:39:53.448657"
$ws.Range("F6").Value = "2021-10-05 13:39:53.448660"
$ws.Range("F7").Value = "2021-10-05 13:39:53.448663"
$ws.Range("F8").Value = "2021-10-05 13:39:53.448665"
$ws.Range("F9").Value = "2021-10-05 13:39:53.448669"
$ws.Range("F10").Value = "2021-10-05 13:39:53.448672"
$ws.Range("F11").Value = "2021-10-05 13:39:53.448675"
$ws.Range("F12").Value = "2021-10-05 13:39:53.448678"
$ws.Range("F13").Value = "2021-10-05 13:39:53.448680"
$ws.Range("F14").Value = "2021-10-05 13:39:53.448683"
$ws.Range("F15").Value = "2021-10-05 13:39:53.448685"
$ws.Range("F16").Value = "2021-10-05 13:39:53.448688"
$ws.Range("F17").Value = "2021-10-05 13:39:53.448691"
$ws.Range("F18").Value = "2021-10-05 13:39:53.448694"
$ws.Range("F19").Value = "2021-10-05 13:39:53.448696"
$ws.Range("F20").Value = "2021-10-05 13:39:53.448699"
$ws.Range("F21").Value = "2021-10-05 13:39:53.448701"
$ws.Range("F22").Value = "2021-10-05 13:39:53.448704"
$ws.Range("F23").Value = "2021-10-05 13:39:53.448706"
$ws.Range("F24").Value = "2021-10-05 13:39:53.448709"
$ws.Range("F25").Value = "2021-10-05 13:39:53.448712"
$ws.Range("F26").Value = "2021-10-05 13:39:53.448714"
$ws.Range("F27").Value = "2021-10-05 13:39:53.448717"
$ws.Range("F28").Value = "2021-10-05 13:39:53.448720"
$ws.Range("F29").Value = "2021-10-05 13:39:53.448722"
$ws.Range("F30").Value = "2021-10-05 13:39:53.448725"
$ws.Range("F31").Value = "2021-10-05 13:39:53.448727"
$ws.Range("F32").Value = "2021-10-05 13:39:53.448730"
$ws.Range("F33").Value = "2021-10-05 13:39:53.448732"
$ws.Range("F34").Value = "2021-10-05 13:39:53.448735"
$ws.Range("F35").Value = "2021-10-05 13:39:53.448738"
$ws.Range("F36").Value = "2021-10-05 13:39:53.448740"
$ws.Range("F37").Value = "2021-10-05 13:39:53.448743"
$ws.Range("F38").Value = "2021-10-05 13:39:53.448745"
$ws.Range("F39").Value = "2021-10-05 13:39:53.448748"
$ws.Range("F40").Value = "2021-10-05 13:39:53.448751"
$ws.Range("F41").Value = "2021-10-05 13:39:53.448753"
$ws.Range("F42").Value = "2021-10-05 13:39:53.448756"
$ws.Range("F43").Value = "2021-10-05 13:39:53.448759"
$ws.Range("F44").Value = "2021-10-05 13:39:53.448762"
$ws.Range("F45").Value = "2021-10-05 13:39:53.448764"
$ws.Range("F46").Value = "2021-10-05 13:39:53.448767"
$ws.Range("F47").Value = "2021-10-05 13:39:53.448769"
$ws.Range("F48").Value = "2021-10-05 13:39:53.448772"
$ws.Range("F49").Value = "2021-10-05 13:39:53.448774"
$ws.Range("F50").Value = "2021-10-05 13:39:53.448777"
$ws.Range("F51").Value = "2021-10-05 13:39:53.448779"
$ws.Range("F52").Value = "2021-10-05 13:39:53.448782"
$ws.Range("F53").Value = "2021-10-05 13:39:53.448784"
$ws.Range("F54").Value = "2021-10-05 13:39:53.448787"
$ws.Range("F55").Value = "2021-10-05 13:39:53.448790"
$ws.Range("F56").Value = "2021-10-05 13:39:53.448792"
$ws.Range("F57").Value = "2021-10-05 13:39:53.448795"
$ws.Range("F58").Value = "2021-10-05 13:39:53.448797"
$ws.Range("F59").Value = "2021-10-05 13:39:53.448800"
$ws.Range("F60").Value = "2021-10-05 13:39:53.448802"
$ws.Range("F61").Value = "2021-10-05 13:39:53.448805"
$ws.Range("F62").Value = "2021-10-05 13:39:53.448807"
$ws.Range("F63").Value = "2021-10-05 13:39:53.448810"
$ws.Range("F64").Value = "2021-10-05 13:39:53.448812"
$ws.Range("F65").Value = "2021-10-05 13:39:53.448815"
$ws.Range("F66").Value = "2021-10-05 13:39:53.448818"
$ws.Range("F67").Value = "2021-10-05 13:39:53.448821"
$ws.Range("F68").Value = "2021-10-05 13:39:53.448824"
$ws.Range("F69").Value = "2021-10-05 13:39:53.448826"
$ws.Range("F70").Value = "2021-10-05 13:39:53.448829"
$ws.Range("F71").Value = "2021-10-05 13:39:53.448832"
$ws.Range("F72").Value = "2021-10-05 13:39:53.448834"
$ws.Range("F73").Value = "2021-10-05 13:39:53.448837"
$ws.Range("F74").Value = "2021-10-05 13:39:53.448839"
$ws.Range("F75").Value = "2021-10-05 13:39:53.448842"
$ws.Range("F76").Value = "2021-10-05 13:39:53.448845"
$ws.Range("F77").Value = "2021-10-05 13:39:53.448847"
$ws.Range("F78").Value = "2021-10-05 13:39:53.448852"
$ws.Range("F79").Value = "2021-10-05 13:39:53.448855"
$ws.Range("F80").Value = "2021-10-05 13:39:53.448857"
$ws.Range("F81").Value = "2021-10-05 13:39:53.448860"
$ws.Range("F82").Value = "2021-10-05 13:39:53.448863"
$ws.Range("F83").Value = "2021-10-05 13:39:53.448865"
$ws.Range("F84").Value = "2021-10-05 13:39:53.448868"
$ws.Range("F85").Value = "2021-10-05 13:39:53.448871"
$ws.Range("F86").Value = "2021-10-05 13:39:53.448873"
$ws.Range("F87").Value = "2021-10-05 13:39:53.448876"
$ws.Range("F88").Value = "2021-10-05 13:39:53.448879"
$ws.Range("F89").Value = "2021-10-05 13:39:53.448881"
$ws.Range("F90").Value = "2021-10-05 13:39:53.448884"
